# Insert a new data row at row 24 (pushing existing rows 24-50 down to 25-51)
# and populate it with the new Damasco record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value = "Maule"
$ws.Cells.Item(24, 4).Value = 44915
$ws.Cells.Item(24, 5).Value = 7
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100103
$ws.Cells.Item(24, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(24, 9).Value = 100103003
$ws.Cells.Item(24, 10).Value = "Damasco"
$ws.Cells.Item(24, 11).Value = "Dina"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 270
$ws.Cells.Item(24, 14).Value = 18000
$ws.Cells.Item(24, 15).Value = 18000
$ws.Cells.Item(24, 16).Value = 18000
$ws.Cells.Item(24, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(24, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 19).Value = 1000
$ws.Cells.Item(24, 20).Value = 18
